# Localize the workbook: translate sheet names from Russian to English.
# Excel automatically rewrites any formulas that reference these sheets
# (e.g. ='Продажи'!B6 -> ='Sales'!B6), so no formula text needs to be
# touched by hand.

$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("Сводный").Name      = "Consolidated budget"
$wb.Worksheets.Item("Продажи").Name      = "Sales"
$wb.Worksheets.Item("Производство").Name = "Production"
$wb.Worksheets.Item("Логистика").Name    = "Logistics"
$wb.Worksheets.Item("Прочее").Name       = "Misc"

# Restore each sheet's active-cell selection (as left by the editor).
$wb.Worksheets.Item("Consolidated budget").Activate()
$wb.Worksheets.Item("Consolidated budget").Range("E10").Select()

$wb.Worksheets.Item("Sales").Activate()
$wb.Worksheets.Item("Sales").Range("C34").Select()

$wb.Worksheets.Item("Production").Activate()
$wb.Worksheets.Item("Production").Range("E34").Select()

$wb.Worksheets.Item("Logistics").Activate()
$wb.Worksheets.Item("Logistics").Range("F33").Select()

$wb.Worksheets.Item("Misc").Activate()
$wb.Worksheets.Item("Misc").Range("H34").Select()

# Leave focus on the consolidated-budget tab (tabSelected="1" in source).
$wb.Worksheets.Item("Consolidated budget").Activate()
$wb.Worksheets.Item("Consolidated budget").Range("E10").Select()
